# Apply revisions to the "Antibodies" example workbook:
#  - Antibodies sheet: replace several antibody-name cells with
#    well-name-like placeholders, fix some host/isotype values, and drop
#    the now-unused last data row (row 11).
#  - Terminology sheet: drop the trailing "kappa"/"lambda" rows (16-17)
#    that are no longer offered, and fix up the Isotype list validation
#    on the Antibodies sheet to match the new, shorter range.

$wb = $excel.ActiveWorkbook
$ab = $wb.Worksheets.Item("Antibodies")
$term = $wb.Worksheets.Item("Terminology")

# --- Antibodies sheet: update individual cell values (rows 2-10) -----
$ab.Range("A2").Value = "A6"

$ab.Range("A3").Value = "B12"
$ab.Range("B3").Value = "Mus musculus"

$ab.Range("A4").Value = ""
$ab.Range("C4").Value = "IgD"

$ab.Range("A5").Value = "C2"
$ab.Range("B5").Value = "Mus musculus"
$ab.Range("C5").Value = "IgG"

$ab.Range("A6").Value = "C3"
$ab.Range("B6").Value = "Homo sapiens"
$ab.Range("C6").Value = "IgG2a"

$ab.Range("A7").Value = "C6"
$ab.Range("B7").Value = ""
$ab.Range("C7").Value = "Ig1"

$ab.Range("A8").Value = "D12"
$ab.Range("B8").Value = "Homo sapiens"

$ab.Range("A9").Value = "E8"
$ab.Range("B9").Value = "Mus musclus"
$ab.Range("C9").Value = "Igm"

$ab.Range("A10").Value = "C3"
$ab.Range("C10").Value = "IgG2a"

# Row 11 ("Acme mAb 10" / "Mus musculus" / "Ig") is no longer needed.
$ab.Rows.Item(11).Delete()

# --- Terminology sheet: drop the last two rows (kappa / lambda) ------
$term.Rows.Item(17).Delete()
$term.Rows.Item(16).Delete()

# --- Restore/fix the data validations on the Antibodies sheet. -------
# Deleting row 11 above shifts existing validation sqrefs, so recreate
# both validations explicitly with their correct extents: Host's list
# is unchanged, Isotype's list now points at the shorter Terminology
# range (B2:B15 instead of B2:B17).
$rb = $ab.Range("B2:B100")
$rb.Validation.Delete()
$rb.Validation.Add(3, 1, 1, "=Terminology!`$A`$2:`$A`$4")
$rb.Validation.IgnoreBlank = $true
$rb.Validation.InCellDropdown = $true
$rb.Validation.ShowInput = $true
$rb.Validation.ShowError = $true

$rc = $ab.Range("C2:C100")
$rc.Validation.Delete()
$rc.Validation.Add(3, 1, 1, "=Terminology!`$B`$2:`$B`$15")
$rc.Validation.IgnoreBlank = $true
$rc.Validation.InCellDropdown = $true
$rc.Validation.ShowInput = $true
$rc.Validation.ShowError = $true
